$wb = $excel.ActiveWorkbook

# Rename the two worksheets to unify the DataNode/DataTable/Entity naming
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "DataNode_1"
$ws2.Name = "DataNode_2"

# Adjust header row heights (row 1) on both sheets
$ws1.Rows.Item(1).RowHeight = 27
$ws2.Rows.Item(1).RowHeight = 27

# Adjust the field-description row heights (row 8) on both sheets
$ws1.Rows.Item(8).RowHeight = 40.5
$ws2.Rows.Item(8).RowHeight = 67.5

# Make the second sheet the active tab
$ws2.Activate()
